$d = $word.ActiveDocument

# NOTE: We deliberately search with Find.Execute (no Replacement arg) and
# then assign the narrowed-down found Range's .Text property directly,
# instead of passing a ReplaceWith string to Find.Execute. The latter
# silently "smart quotes" straight apostrophes ' into curly ones while
# going through Word's replace-text pipeline; direct Range.Text assignment
# preserves the literal characters we supply (matching the source diff,
# which uses straight apostrophes throughout).

function Replace-ParagraphText($oldText, $newText) {
    $range = $d.Content
    $found = $range.Find.Execute($oldText)
    if (-not $found) {
        throw "Could not find expected text: $oldText"
    }
    $range.Text = $newText
}

Replace-ParagraphText `
    'The economic landscape of San Diego County has been marked by a dynamic evolution in its Gross Regional Product (GRP) over recent years. This analysis delves into the GRP data from 2019 to 2023, providing insights into the economic vitality of the region, with comparisons to state and national levels.' `
    'The economic landscape of San Diego County has been marked by a dynamic evolution in its Gross Regional Product (GRP) over recent years. This analysis delves into the GRP data from 2019 to 2023, offering insights into the economic vitality of the region, with comparisons to state and national levels.'

Replace-ParagraphText `
    'In 2019, San Diego County''s GRP stood at approximately $244.28 billion, with a per capita GRP of $73,347. This figure slightly increased in 2020 to $244.82 billion, despite a decrease in population, resulting in a per capita GRP of $74,278. The year 2021 marked a significant upturn, with the GRP rising to $268.87 billion and the per capita GRP reaching $82,100. This upward trajectory continued into 2022, with the GRP climbing to $296.68 billion and a per capita GRP of $90,557. By 2023, the GRP further increased to $308.71 billion, with a per capita GRP of $94,916.' `
    'In 2019, San Diego County''s GRP stood at approximately $244.28 billion, with a population of 3,330,458, resulting in a per capita GRP of $73,347. This figure was slightly below California''s per capita GRP of $75,789 but significantly higher than the national average of $63,754. The following year, 2020, saw a modest increase in San Diego''s GRP to $244.82 billion, despite a slight population decline to 3,296,045. This led to a per capita GRP of $74,278, reflecting resilience amid broader economic challenges. California''s per capita GRP dipped slightly to $74,964, while the national figure decreased to $62,157, indicating a more pronounced impact of economic disruptions at the state and national levels.'

Replace-ParagraphText `
    'Comparatively, California''s GRP also demonstrated growth over the same period. In 2019, the state''s GRP was approximately $2.99 trillion, with a per capita GRP of $75,789. This figure saw a slight decline in 2020 to $2.96 trillion, with a per capita GRP of $74,964. However, by 2021, the state''s GRP had surged to $3.31 trillion, with a per capita GRP of $84,587. The growth continued in 2022, reaching $3.54 trillion and a per capita GRP of $90,636, and further to $3.65 trillion in 2023, with a per capita GRP of $93,800.' `
    'The year 2021 marked a significant upturn for San Diego County, with the GRP rising to $268.87 billion. The population continued to decrease slightly to 3,274,954, but the per capita GRP surged to $82,100. This growth outpaced both California''s per capita GRP of $84,587 and the national average of $68,858, highlighting San Diego''s robust economic recovery and expansion. In 2022, the county''s GRP further increased to $296.68 billion, with a stable population of 3,276,208. The per capita GRP reached $90,557, continuing to exceed the state average of $90,636 and the national figure of $74,889. This trend underscores San Diego''s sustained economic momentum and its ability to leverage its diverse economic base.'

Replace-ParagraphText `
    'On a national scale, the United States'' GRP in 2019 was approximately $20.93 trillion, with a per capita GRP of $63,754. This figure decreased in 2020 to $20.61 trillion, with a per capita GRP of $62,157. By 2021, the national GRP had rebounded to $22.86 trillion, with a per capita GRP of $68,858. The upward trend persisted in 2022, with the GRP reaching $24.96 trillion and a per capita GRP of $74,889, and further to $25.96 trillion in 2023, with a per capita GRP of $77,366.' `
    'By 2023, San Diego County''s GRP had grown to $308.71 billion, with a population of 3,252,468, resulting in a per capita GRP of $94,916. This figure remained competitive with California''s per capita GRP of $93,800 and significantly above the national average of $77,366. The data reflects San Diego''s economic resilience and adaptability, driven by key sectors such as technology, tourism, and defense, which have continued to thrive and contribute to the region''s economic prosperity.'

Replace-ParagraphText `
    'The data reveals that San Diego County''s economic growth has been robust, outpacing both state and national averages in terms of per capita GRP. This growth reflects the region''s resilience and adaptability in the face of economic challenges, underscoring its role as a significant contributor to the broader economic landscape. As San Diego County continues to expand its economic footprint, it remains a vital player in the state''s and nation''s economic narratives.' `
    'Overall, San Diego County''s economic performance over this period illustrates a pattern of growth and resilience, with its GRP consistently rising and its per capita figures outperforming national averages. This trajectory not only highlights the county''s economic strengths but also positions it as a vital contributor to California''s overall economic health. As San Diego continues to navigate the complexities of the global economy, its strategic focus on innovation and diversification will likely sustain its economic vitality in the years to come.'

Write-Host "Done: all 5 paragraph segments replaced."
